$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(68).Insert()
# Now copy formatting from row 69 (which now holds old row 68's content/style) into row 68
$ws.Rows.Item(69).Copy()
$ws.Rows.Item(68).PasteSpecial(-4122)  # xlPasteFormats = -4122
Write-Host "done"
